# "Generate Report for Handoff" - refresh the localization-status report:
#   - Status moves from "Handed back: in sync with en-US" to "Ready for handoff"
#   - the handoff timestamps are bumped to the new generation time
#   - the (now shorter) status text lets Excel narrow the status columns

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- Timestamps bumped to the new handoff generation run ---
$overview.Range("G2").Value = "2016-08-18 11:02:45"
$dede.Range("H2").Value     = "2016-08-18 11:02:45"
$zhcn.Range("H2").Value     = "2016-08-18 11:02:41"

# --- Status columns narrow now that the status text is shorter ---
$overview.Columns.Item(5).ColumnWidth = 16.3333333333333
$overview.Columns.Item(6).ColumnWidth = 16.3333333333333
$zhcn.Columns.Item(3).ColumnWidth     = 16.3333333333333
$dede.Columns.Item(3).ColumnWidth     = 16.3333333333333

Write-Host "Applied handoff report refresh"
